# Append a new paragraph "Let the learning continue" (with the same
# run-splitting / gramStart-gramEnd proofing markers the source document
# uses around "together") followed by a new empty paragraph, mirroring
# the diff that adds this content right after the existing paragraph and
# before the sectPr.

$d = $word.ActiveDocument

# Collapse to the very end of the document body content, then insert the
# new paragraphs there (i.e. right before the final section properties).
$r = $d.Content
$r.Collapse(0)

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Let the learning </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>cont</w:t></w:r><w:r><w:t>inue</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$r.InsertXML($xml)
